# Apply TABLE_CONFIG column width specifications.
#
# Slide 5 - "Timeline & Milestones" table: resize the 4 grid columns.
# Slide 8 - "Investment Summary" table: resize the 7 grid columns.
#
# PowerPoint's Table.Columns(i).Width property is expressed in points;
# 1 point = 12700 EMU. The target widths below (in EMU) come straight
# from the desired <a:gridCol w="..."/> values, converted to points.

$p = $ppt.ActivePresentation

# ---- Slide 5: Timeline & Milestones table -----------------------------
$s5 = $p.Slides.Item(5)
$tbl5 = $s5.Shapes.Item(3).Table

$widths5 = @(871093, 2177733, 1306639, 4355466)
for ($i = 1; $i -le $widths5.Count; $i++) {
    $tbl5.Columns.Item($i).Width = $widths5[$i - 1] / 12700.0
}

# ---- Slide 8: Investment Summary table --------------------------------
$s8 = $p.Slides.Item(8)
$tbl8 = $s8.Shapes.Item(3).Table

$widths8 = @(1742186, 871093, 1829295, 1219530, 958202, 958202, 1132421)
for ($i = 1; $i -le $widths8.Count; $i++) {
    $tbl8.Columns.Item($i).Width = $widths8[$i - 1] / 12700.0
}
